$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear row 3 values first (B3 cell emptied entirely, A3 keeps its style but no value)
$ws.Range("A3").Value = $null
$ws.Range("B3").Value = $null

# Update row 2 values (B2 before A2 to match shared-string ordering)
$ws.Range("B2").Value = "Rakesh11"
$ws.Range("A2").Value = "AutomationCategory99"

# Widen column A (closest achievable value to 25.140625 given pixel rounding)
$ws.Columns.Item(1).ColumnWidth = 24.33

# Update selection to A3
$ws.Range("A3").Select()
